$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sharedStrings table gains two new entries ("line7","line8") inserted right
# after "line6" / right before "extr1". This shifts the text shown by the
# existing rows 8-15 (their underlying shared-string slot keeps the same
# position in the table, but that slot now holds different text), and two
# brand new rows (16,17) are appended at the end of the sheet.
#
# Resulting name column (rows 8-17): line7, line8, extr1, extr2, extr3, extr4,
# extr5, extr6, extr7, extr8 -- with updated from_bus/to_bus/in_service values.

$rows = @(
    @{ Row = 8;  Name = "line7"; From = 14; To = 11; InService = $true  },
    @{ Row = 9;  Name = "line8"; From = 16; To = 9;  InService = $true  },
    @{ Row = 10; Name = "extr1"; From = 5;  To = 12; InService = $true  },
    @{ Row = 11; Name = "extr2"; From = 5;  To = 9;  InService = $true  },
    @{ Row = 12; Name = "extr3"; From = 10; To = 11; InService = $true  },
    @{ Row = 13; Name = "extr4"; From = 7;  To = 8;  InService = $false },
    @{ Row = 14; Name = "extr5"; From = 9;  To = 11; InService = $false },
    @{ Row = 15; Name = "extr6"; From = 7;  To = 11; InService = $false },
    @{ Row = 16; Name = "extr7"; From = 5;  To = 7;  InService = $true  },
    @{ Row = 17; Name = "extr8"; From = 8;  To = 5;  InService = $true  }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.Name
    $ws.Cells.Item($row, 3).Value = $r.From
    $ws.Cells.Item($row, 4).Value = $r.To
    $ws.Cells.Item($row, 5).Value = $r.InService
}

# New rows 16/17 also need the "id" column (A) populated, matching the style
# used by the other column-A cells (bold, bordered, centered) -- copy the
# formatting from an existing column-A cell instead of rebuilding it property
# by property (avoids generating unused intermediate cell styles).
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(17, 1).Value = 15

$ws.Cells.Item(15, 1).Copy() | Out-Null
$ws.Cells.Item(16, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17, 1).PasteSpecial(-4122) | Out-Null
